# Rename the worksheet to reflect the new update date (02-12-2025 -> 05-12-2025).
# Excel automatically updates any defined names / formulas that reference the
# sheet by name when it is renamed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Opdateret d. 05-12-2025"
